# Applies the data refresh captured in the commit "Updated cryptos list on
# Sun Oct 22 17:37:15 UTC 2023 with GitHub Actions" to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) stores its values as plain text in the source feed, even when
# a value happens to look like a plain number (e.g. "215.12", "9.45", "1.00").
# Prefixing the new text with a leading apostrophe forces Excel to keep storing
# it as text instead of silently re-interpreting it as a numeric value, which
# would change its type/rounding (e.g. "1.00" -> 1, "0.0610" -> 0.061).
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
}

Set-TextValue $ws.Range("D2") '29.906.97'
$ws.Range("E2").Value = '  +0.73%  '

Set-TextValue $ws.Range("D3") '1.632.87'
$ws.Range("E3").Value = '  +1.07%  '

$ws.Range("E4").Value = '  +0.78%  '

Set-TextValue $ws.Range("D5") '215.12'
$ws.Range("E5").Value = '  +1.16%  '

Set-TextValue $ws.Range("D6") '0.521'
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("E7").Value = '  +0.79%  '

Set-TextValue $ws.Range("D8") '28.68'
$ws.Range("E8").Value = '  -1.84%  '

$ws.Range("E9").Value = '  +0.52%  '

Set-TextValue $ws.Range("D10") '0.0610'
$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("E11").Value = '  -0.92%  '

Set-TextValue $ws.Range("D12") '1.864.94'
$ws.Range("E12").Value = '  +1.04%  '

Set-TextValue $ws.Range("D13") '1.633.56'
$ws.Range("E13").Value = '  +1.43%  '

Set-TextValue $ws.Range("D14") '0.577'
$ws.Range("E14").Value = '  +1.24%  '

Set-TextValue $ws.Range("D15") '9.45'
$ws.Range("E15").Value = '  +4.53%  '

Set-TextValue $ws.Range("D16") '29.908.31'
$ws.Range("E16").Value = '  +0.79%  '

$ws.Range("E17").Value = '  -1.88%  '

Set-TextValue $ws.Range("D18") '65.46'
$ws.Range("E18").Value = '  +1.99%  '

Set-TextValue $ws.Range("D19") '241.07'
$ws.Range("E19").Value = '  -0.28%  '

Set-TextValue $ws.Range("D20") '0.0₃0704'
$ws.Range("E20").Value = '  -1.16%  '

Set-TextValue $ws.Range("D21") '1.00'
$ws.Range("E21").Value = '  +0.59%  '

Set-TextValue $ws.Range("D22") '9.86'
$ws.Range("E22").Value = '  +1.69%  '

Set-TextValue $ws.Range("D23") '4.14'
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("E24").Value = '  +2.86%  '

Set-TextValue $ws.Range("D25") '157.68'
$ws.Range("E25").Value = '  +0.80%  '

Set-TextValue $ws.Range("D26") '15.51'
$ws.Range("E26").Value = '  -0.97%  '

$ws.Range("E27").Value = '  -1.26%  '

Set-TextValue $ws.Range("D28") '6.63'
$ws.Range("E28").Value = '  +0.38%  '

Set-TextValue $ws.Range("D29") '0.999'
$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  +1.04%  '

Set-TextValue $ws.Range("D32") '3.39'
$ws.Range("E32").Value = '  +1.79%  '

$ws.Range("E33").Value = '  -1.00%  '

Set-TextValue $ws.Range("D34") '1.425.95'
$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("E35").Value = '  +3.48%  '

$ws.Range("E36").Value = '  -2.52%  '

$ws.Range("E37").Value = '  -2.88%  '

$ws.Range("E38").Value = '  +0.44%  '

$ws.Range("E39").Value = '  +0.20%  '

Set-TextValue $ws.Range("D40") '75.53'
$ws.Range("E40").Value = '  +8.16%  '

Set-TextValue $ws.Range("D41") '0.557'
$ws.Range("E41").Value = '  -0.24%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D42") '0.836'
$ws.Range("E42").Value = '  +0.54%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D43") '1.99'
$ws.Range("E43").Value = '  +1.40%  '

Set-TextValue $ws.Range("D44") '0.0501'
$ws.Range("E44").Value = '  -0.89%  '

$ws.Range("E45").Value = '  +0.76%  '

$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("E47").Value = '  -1.78%  '

Set-TextValue $ws.Range("D48") '1.772.23'
$ws.Range("E48").Value = '  +0.99%  '

Set-TextValue $ws.Range("D49") '48.44'
$ws.Range("E49").Value = '  -10.18%  '

Set-TextValue $ws.Range("D50") '92.48'
$ws.Range("E50").Value = '  +5.18%  '

$ws.Range("E51").Value = '  +3.87%  '
